# "Agregando soporte para apoyos inclinados"
# Adds a new "rotación" (rotation) column to the restric (supports) sheet,
# so inclined supports can be modeled with their rotation angle in degrees.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("restric")

# --- New column D: rotación ----------------------------------------------

# Header cell, styled like the existing A1/C1 headers.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "rotación"
$ws.Range("D1").AddComment("grados")

# Nodes 1-3 (rows 2-4) have no rotational restraint defined yet -> 0.
$ws.Range("D2:D4").Value = 0

# Nodes with rotational data pending (rows 5-8) default to #N/A, matching
# the style (General format + yellow highlight) already used in column C.
$ws.Range("D5").Formula = "=NA()"
$ws.Range("D6").Formula = "=NA()"
$ws.Range("D7").Formula = "=NA()"
$ws.Range("D8").Formula = "=NA()"
$ws.Range("C5").Copy()
$ws.Range("D5:D8").PasteSpecial(-4122)

# --- Make "restric" the active sheet/cell, as in the saved workbook ------

$ws.Activate()
$ws.Range("F9").Select()
